# Update countries & provincias Spain
# Applies the data refresh captured in the commit: a few country case
# counts were updated, Venezuela overtook Malta in the ranking, the
# Seychelles/Groenlandia/Montserrat trio and the San Bartolome/Sahara
# Occidental pair were re-sorted, and the "last updated" timestamp moved
# from 02:05 to 02:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 02:35"

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Range("B4").Value = 1550083
$ws.Range("C4").Value = 22419
$ws.Range("D4").Value = 355677
$ws.Range("E4").Value = 1102430
$ws.Range("G4").Value = 998
$ws.Range("H4").Value = 91976

# --- Venezuela overtakes Malta (rows 124-125 swap, Venezuela gets fresh data) ---
$ws.Range("A124").Value = "Venezuela"
$ws.Range("B124").Value = 618
$ws.Range("C124").Value = 77
$ws.Range("D124").Value = 241
$ws.Range("E124").Value = 367
$ws.Range("H124").Value = 10

$ws.Range("A125").Value = "Malta"
$ws.Range("B125").Value = 558
$ws.Range("C125").Value = 5
$ws.Range("D125").Value = 456
$ws.Range("E125").Value = 96
$ws.Range("H125").Value = 6

# --- Bahamas (row 168): D/E swap ---
$ws.Range("D168").Value = 43
$ws.Range("E168").Value = 42

# --- Seychelles / Groenlandia / Montserrat re-sort (rows 209-211) ---
$ws.Range("A209").Value = "Seychelles"

$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- San Bartolome / Sahara Occidental swap (rows 214 & 216) ---
$ws.Range("A214").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"
